$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats/styles) from row 8 into row 9 first,
# so the new row's A9/G9 cells reuse the existing date/bool styles instead
# of creating brand new ones. -4122 == xlPasteFormats.
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)

$ws.Cells.Item(9, 1).Value = 42654.743703703702
$ws.Cells.Item(9, 2).Value = $true
$ws.Cells.Item(9, 3).Value = 10071.69
$ws.Cells.Item(9, 4).Value = 10071.19
$ws.Cells.Item(9, 5).Value = 75.5
$ws.Cells.Item(9, 6).Value = 75.489998
$ws.Cells.Item(9, 7).Value = $true
$ws.Cells.Item(9, 8).Value = -0.01
$ws.Cells.Item(9, 9).Value = $false

# Re-run the "best fit" column sizing now that the new row has changed the
# widest content in each column (mirrors hitting AutoFit on the sheet).
$ws.Columns.Item(1).ColumnWidth = 14.541666666666666
$ws.Columns.Item(2).ColumnWidth = 7.416666666666667
$ws.Columns.Item(3).ColumnWidth = 8.041666666666666
$ws.Columns.Item(4).ColumnWidth = 10.416666666666666
$ws.Columns.Item(5).ColumnWidth = 9.041666666666666
$ws.Columns.Item(6).ColumnWidth = 9.041666666666666
$ws.Columns.Item(7).ColumnWidth = 9.541666666666666
$ws.Columns.Item(8).ColumnWidth = 13.791666666666666
$ws.Columns.Item(9).ColumnWidth = 11.041666666666666
